$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("patient_id", "ethnicity", "race"),
    @(11332, "Hispanic or Latino", "Black or African American"),
    @(20202, "Hispanic or Latino", "White"),
    @(26287, "Not Hispanic or Latino", "Asian"),
    @(32851, "Not Hispanic or Latino", "Black or African American"),
    @(34899, "Hispanic or Latino", "White"),
    @(35383, "Not Hispanic or Latino", "Other Race"),
    @(36963, "Not Hispanic or Latino", "Asian"),
    @(44221, "Not Hispanic or Latino", "Black or African American"),
    @(45421, "Not Hispanic or Latino", "White"),
    @(49143, "Hispanic or Latino", "American Indian or Alaska Native"),
    @(56517, "Not Hispanic or Latino", "Other Race"),
    @(58898, "Not Hispanic or Latino", "Asian"),
    @(64017, "Not Hispanic or Latino", "Black or African American"),
    @(64607, "Not Hispanic or Latino", "Asian"),
    @(76051, "Hispanic or Latino", "Other Race"),
    @(90185, "Hispanic or Latino", "Native Hawaiian or Other Pacific Islander"),
    @(92210, "Not Hispanic or Latino", "American Indian or Alaska Native")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 1 -and $row -gt 1) {
            # patient_id values look numeric but must remain stored as text
            $cell.NumberFormat = "@"
            $cell.Value = [string]$rowData[$j]
            $cell.NumberFormat = "General"
        } else {
            $cell.Value = [string]$rowData[$j]
        }
    }
}
